$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arbeitsmatrix")

# --- Row 111 -------------------------------------------------------------
# Copy formatting from row 110 (same visual "block" style, font color s=32)
# for the cells that need that look, then overwrite values/formulas.
$ws.Cells.Item(110,1).Copy($ws.Cells.Item(111,1)) | Out-Null    # A111 s=32
$ws.Cells.Item(110,2).Copy($ws.Cells.Item(111,2)) | Out-Null    # B111 s=32
$ws.Cells.Item(110,3).Copy($ws.Cells.Item(111,3)) | Out-Null    # C111 s=32
$ws.Cells.Item(110,5).Copy($ws.Cells.Item(111,5)) | Out-Null    # E111 s=32
$ws.Cells.Item(110,9).Copy($ws.Cells.Item(111,9)) | Out-Null    # I111 s=15
$ws.Cells.Item(110,10).Copy($ws.Cells.Item(111,10)) | Out-Null  # J111 s=21
$ws.Cells.Item(110,11).Copy($ws.Cells.Item(111,11)) | Out-Null  # K111 s=22

$ws.Cells.Item(111,1).Value = 22
$ws.Cells.Item(111,2).Value = "Interface Design"
$ws.Cells.Item(111,3).Value = "MockUps"
$ws.Cells.Item(111,4).Value = "[FEATURE]"
$ws.Cells.Item(111,5).Value = "MockUps Rezept Filter"
$ws.Cells.Item(111,6).Value = 44380
$ws.Cells.Item(111,7).Value = 44359
$ws.Cells.Item(111,10).Value = 0.58333333333333337
$ws.Cells.Item(111,11).Value = 0.66666666666666663

# --- Row 112 ---------------------------------------------------------------
# Row 112 keeps the default (unstyled) look for A/B/C/E, matching rows
# such as 107/108 that have no explicit "s" attribute on those columns.
$ws.Cells.Item(107,1).Copy($ws.Cells.Item(112,1)) | Out-Null    # A112 (no style)
$ws.Cells.Item(107,2).Copy($ws.Cells.Item(112,2)) | Out-Null    # B112 (no style)
$ws.Cells.Item(107,3).Copy($ws.Cells.Item(112,3)) | Out-Null    # C112 (no style)
$ws.Cells.Item(107,5).Copy($ws.Cells.Item(112,5)) | Out-Null    # E112 (no style)
$ws.Cells.Item(110,9).Copy($ws.Cells.Item(112,9)) | Out-Null    # I112 s=15
$ws.Cells.Item(110,10).Copy($ws.Cells.Item(112,10)) | Out-Null  # J112 s=21
$ws.Cells.Item(110,11).Copy($ws.Cells.Item(112,11)) | Out-Null  # K112 s=22

$ws.Cells.Item(112,1).Value = 18
$ws.Cells.Item(112,2).Value = "Konzeptuelles Design"
$ws.Cells.Item(112,3).Value = "Content Map"
$ws.Cells.Item(112,4).Value = "[FEATURE]"
$ws.Cells.Item(112,5).Value = "Filterarten bestimmen und konzipieren"
$ws.Cells.Item(112,6).Value = 44380
$ws.Cells.Item(112,7).Value = 44359
$ws.Cells.Item(112,10).Value = 0.66666666666666663
$ws.Cells.Item(112,11).Value = 0.70833333333333337

# --- Formulas (written as two separate shared-formula groups, matching
#     the original I105:I110 group plus a new I111:I112 group) -----------
$ws.Range("I105:I110").Formula = "=ROUNDUP(((SUM(K105-J105)*24*60/60)/0.25),0)*0.25"
$ws.Range("I111:I112").Formula = "=ROUNDUP(((SUM(K111-J111)*24*60/60)/0.25),0)*0.25"

# --- Sheet view state ------------------------------------------------------
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 100
$win.ScrollColumn = 1
$ws.Range("C103").Select() | Out-Null
